$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.761.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.948.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4818"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2960"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06817"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "112.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.954.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.559"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07645"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6898"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "298.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.778.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007704"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.641"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.204.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.587"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.709"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1085"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.432"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.588"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.384"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05061"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7729"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.165"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02079"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "111.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4475"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8763"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.923"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.397"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.479"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2551"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
